# Apply "updated contents and readme" edit to the codebook worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 3: the "combine male/female" caveat note is no longer needed - clear it.
$ws.Range("E3").ClearContents()

# 2. Row 5: clarify the racial-composition description is expressed as a percentage.
$ws.Range("D5").Value = "Racial composition by percentage - black nonhisp, white nonhisp, asian nonhisp, hispanic, two or more"

# 3. Row 9: replace the placeholder "S" with the real description.
$ws.Range("D9").Value = "Percent below poverty level"

# 4. Row 13: fix capitalization of the note.
$ws.Range("E13").Value = "Not available 2009-2010"

# 5. Row 17: fill in the new population_density entry that was only half-populated.
$ws.Range("A17").Value = "NA"
$ws.Range("B17").Value = "NA"
$ws.Range("E17").Value = "TODO: calculate after spatial join"
$ws.Range("F17").Value = "WIP"

# 6. Widen columns D and E to fit the longer text.
$ws.Range("D1").EntireColumn.ColumnWidth = 57.83
$ws.Range("E1").EntireColumn.ColumnWidth = 30.17

# 7. Leave the selection on D10, matching the saved view state.
[void]$ws.Range("D10").Select()
